$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update threshold values (row 2: alpha_distance_range, row 3: beta_distance_range, row 4: ratio_threshold_range)
$ws.Range("B2").Value = 5.3
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 8.3
$ws.Range("B4").Value = 0.75

# Move the active selection from C5 to C4, matching the saved view state
$ws.Range("C4").Select() | Out-Null
